$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": rows 3 and 4 swap identity (ac3fb723 moves up to row 3,
# 1d4eb241 moves down to row 4). New "Ready for handoff" status + new
# timestamp for 1d4eb241 / d27dba92.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A3").Value = "ac3fb723-3934-4251-86ab-b855325dbdf1.md"
$ws1.Range("B3").Value = "e2e\ac3fb723-3934-4251-86ab-b855325dbdf1.md"
$ws1.Range("C3").Value = ".md"
$ws1.Range("D3").Value = ""
$ws1.Range("E3").Value = "Handed back: in sync with en-US"
$ws1.Range("F3").Value = "Handed back: in sync with en-US"
$ws1.Range("G3").Value = "2016-08-17 22:22:56"

$ws1.Range("A4").Value = "1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.md"
$ws1.Range("B4").Value = "e2e\1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.md"
$ws1.Range("C4").Value = ".md"
$ws1.Range("D4").Value = ""
$ws1.Range("E4").Value = "Ready for handoff"
$ws1.Range("F4").Value = "Ready for handoff"
$ws1.Range("G4").Value = "2016-08-17 22:25:35"

$ws1.Range("E5").Value = "Ready for handoff"
$ws1.Range("F5").Value = "Ready for handoff"
$ws1.Range("G5").Value = "2016-08-17 22:25:35"

# Hyperlinks: B3 now points at the "ac3fb723" display text, B4 at "1d4eb241".
# The underlying relationship ids keep their original targets (unchanged in
# the source diff), so we only touch the link TextToDisplay captions here by
# re-creating them in the swapped order while keeping the same targets.
$targetB3 = $ws1.Hyperlinks.Item(2).Address
$targetB4 = $ws1.Hyperlinks.Item(3).Address
$ws1.Hyperlinks.Item(2).Delete()
$ws1.Hyperlinks.Item(2).Delete()
$ws1.Hyperlinks.Add($ws1.Range("B3"), $targetB3, [Type]::Missing, [Type]::Missing, "e2e\ac3fb723-3934-4251-86ab-b855325dbdf1.md")
$ws1.Hyperlinks.Add($ws1.Range("B4"), $targetB4, [Type]::Missing, [Type]::Missing, "e2e\1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn": detail rows for ac3fb723 (row3, now fully populated),
# 1d4eb241 (row4) and d27dba92 (row5) get new statuses / error detail.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Columns.Item(16).ColumnWidth = 40

$ws2.Range("A3").Value = "ac3fb723-3934-4251-86ab-b855325dbdf1.md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("G3").Value = "ac3fb723-3934-4251-86ab-b855325dbdf1.4ea0bbd83bb899adbf7afd34e2211aa107e5fa86.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-17 22:22:50"
$ws2.Range("I3").Value = "ac3fb723-3934-4251-86ab-b855325dbdf1.md"
$ws2.Range("J3").Value = "ac3fb723-3934-4251-86ab-b855325dbdf1.4ea0bbd83bb899adbf7afd34e2211aa107e5fa86.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-17 22:23:51"
$ws2.Range("M3").Value = "True"
$ws2.Range("P3").Value = ""

$ws2.Range("A4").Value = "1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.md"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("G4").Value = "1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.2100c1d4e488aa04321dcb4b2f022cdbb0848c23.zh-cn.xlf"
$ws2.Range("H4").Value = "2016-08-17 22:25:30"
$ws2.Range("I4").Value = "1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.md"
$ws2.Range("J4").Value = "1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.2100c1d4e488aa04321dcb4b2f022cdbb0848c23.zh-cn.xlf"
$ws2.Range("K4").Value = "2016-08-17 22:24:49"
$ws2.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/343f1a65e9f1d91e1d14ab80e26924182a29fc63/e2e/1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d965c32e99a258d9aa65ee2bfa51d1761e05392e/e2e/1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.md."

$ws2.Range("C5").Value = "Ready for handoff"
$ws2.Range("H5").Value = "2016-08-17 22:25:30"
$ws2.Range("K5").Value = "2016-08-17 22:24:49"
$ws2.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/343f1a65e9f1d91e1d14ab80e26924182a29fc63/e2e/d27dba92-0beb-41df-b079-f0b92552fef3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d965c32e99a258d9aa65ee2bfa51d1761e05392e/e2e/d27dba92-0beb-41df-b079-f0b92552fef3.md."

$targetA3_2 = $ws2.Hyperlinks.Item(3).Address
$targetI3_2 = $ws2.Hyperlinks.Item(4).Address
$targetA4_2 = $ws2.Hyperlinks.Item(5).Address
$targetI4_2 = $ws2.Hyperlinks.Item(6).Address
$ws2.Hyperlinks.Item(3).Delete()
$ws2.Hyperlinks.Item(3).Delete()
$ws2.Hyperlinks.Item(3).Delete()
$ws2.Hyperlinks.Item(3).Delete()
$ws2.Hyperlinks.Add($ws2.Range("A3"), $targetA3_2, [Type]::Missing, [Type]::Missing, "ac3fb723-3934-4251-86ab-b855325dbdf1.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), $targetI3_2, [Type]::Missing, [Type]::Missing, "ac3fb723-3934-4251-86ab-b855325dbdf1.md")
$ws2.Hyperlinks.Add($ws2.Range("A4"), $targetA4_2, [Type]::Missing, [Type]::Missing, "1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.md")
$ws2.Hyperlinks.Add($ws2.Range("I4"), $targetI4_2, [Type]::Missing, [Type]::Missing, "1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.md")

# ---------------------------------------------------------------------------
# Sheet "de-de": same pattern as zh-cn.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Columns.Item(16).ColumnWidth = 40

$ws3.Range("A3").Value = "ac3fb723-3934-4251-86ab-b855325dbdf1.md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("G3").Value = "ac3fb723-3934-4251-86ab-b855325dbdf1.4ea0bbd83bb899adbf7afd34e2211aa107e5fa86.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-17 22:22:56"
$ws3.Range("I3").Value = "ac3fb723-3934-4251-86ab-b855325dbdf1.md"
$ws3.Range("J3").Value = "ac3fb723-3934-4251-86ab-b855325dbdf1.4ea0bbd83bb899adbf7afd34e2211aa107e5fa86.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-17 22:23:58"
$ws3.Range("M3").Value = "True"
$ws3.Range("P3").Value = ""

$ws3.Range("A4").Value = "1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.md"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("G4").Value = "1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.2100c1d4e488aa04321dcb4b2f022cdbb0848c23.de-de.xlf"
$ws3.Range("H4").Value = "2016-08-17 22:25:35"
$ws3.Range("I4").Value = "1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.md"
$ws3.Range("J4").Value = "1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.2100c1d4e488aa04321dcb4b2f022cdbb0848c23.de-de.xlf"
$ws3.Range("K4").Value = "2016-08-17 22:24:56"
$ws3.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/343f1a65e9f1d91e1d14ab80e26924182a29fc63/e2e/1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d965c32e99a258d9aa65ee2bfa51d1761e05392e/e2e/1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.md."

$ws3.Range("C5").Value = "Ready for handoff"
$ws3.Range("H5").Value = "2016-08-17 22:25:35"
$ws3.Range("K5").Value = "2016-08-17 22:24:56"
$ws3.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/343f1a65e9f1d91e1d14ab80e26924182a29fc63/e2e/d27dba92-0beb-41df-b079-f0b92552fef3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d965c32e99a258d9aa65ee2bfa51d1761e05392e/e2e/d27dba92-0beb-41df-b079-f0b92552fef3.md."

$targetA3_3 = $ws3.Hyperlinks.Item(3).Address
$targetI3_3 = $ws3.Hyperlinks.Item(4).Address
$targetA4_3 = $ws3.Hyperlinks.Item(5).Address
$targetI4_3 = $ws3.Hyperlinks.Item(6).Address
$ws3.Hyperlinks.Item(3).Delete()
$ws3.Hyperlinks.Item(3).Delete()
$ws3.Hyperlinks.Item(3).Delete()
$ws3.Hyperlinks.Item(3).Delete()
$ws3.Hyperlinks.Add($ws3.Range("A3"), $targetA3_3, [Type]::Missing, [Type]::Missing, "ac3fb723-3934-4251-86ab-b855325dbdf1.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), $targetI3_3, [Type]::Missing, [Type]::Missing, "ac3fb723-3934-4251-86ab-b855325dbdf1.md")
$ws3.Hyperlinks.Add($ws3.Range("A4"), $targetA4_3, [Type]::Missing, [Type]::Missing, "1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.md")
$ws3.Hyperlinks.Add($ws3.Range("I4"), $targetI4_3, [Type]::Missing, [Type]::Missing, "1d4eb241-78ee-4605-9c1f-4a85e9a7c02d.md")
